$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: title paragraph "з дисципліни “..." -> "на тему “..."
# The run "з дисципліни “" becomes "на тему “", with the now-current
# edit point (after the first typed letter "н") marked by the _GoBack
# bookmark, exactly like Word leaves behind after an in-place retype.
# ---------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute("з дисципліни “", $true, $false, $false, $false, $false, $true, 1, $false, "на тему “", 2)

$rng2 = $d.Content.Duplicate
$rng2.Find.Execute("на тему “", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng2.Start

# "_GoBack" marks the point right after the "н" that was typed first.
$editPoint = $start + 1
$bmRange = $d.Range($editPoint, $editPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

# A transient bookmark at the "а тему" / " “" boundary keeps those two
# (identically formatted) spans from being re-flattened into one run;
# removing it again leaves the split in place.
$splitPoint = $start + 7
$splitRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("ZZZ_TEMP_SPLIT", $splitRange)
$d.Bookmarks("ZZZ_TEMP_SPLIT").Delete()

# ---------------------------------------------------------------------
# Edit 2: collapse the three runs that spell out the report title in
# the "на тему:" block into a single run.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "^“Виконання теоретико-множинних операцій реляційної алгебри засобами SQL”") {
        $pr = $p.Range
        $pr.Find.Execute("“Виконання теоретико-множинних операцій реляційної алгебри засобами SQL”", $true, $false, $false, $false, $false, $true, 1, $false, "“Виконання теоретико-множинних операцій реляційної алгебри засобами SQL”", 2)
        break
    }
}

# ---------------------------------------------------------------------
# Edit 3: merge "Перед виконанням..." / " контактів" / ". Результат..."
# into a single run.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "^Перед виконанням завдання") {
        $pr = $p.Range
        $old = "Перед виконанням завдання, потрібно сформувати дві таблиці з однаковими множинами атрибутів. Візьмемо за основу таблицю контактів. Результат збережемо в таблицях "
        $pr.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $old, 2)
        break
    }
}

# ---------------------------------------------------------------------
# Edit 5: merge "Запит на виконання об’єднання" / " Contacts1 і " into
# a single run.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "^Запит на виконання об’єднання Contacts1") {
        $pr = $p.Range
        $old = "Запит на виконання об’єднання Contacts1 і "
        $pr.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $old, 2)
        break
    }
}
